$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1073.7333
$ws.Range("I6").Value = 266.14285
$ws.Range("K6").Value = 798.4285500000001
$ws.Range("M6").Value = -686.4285500000001

$ws.Range("H31").Value = 78.666664
$ws.Range("I31").Value = 78.666664
$ws.Range("K31").Value = 235.999992
$ws.Range("M31").Value = -5.999991999999992

$ws.Range("H62").Value = 5154.6
$ws.Range("I62").Value = 4295
$ws.Range("K62").Value = 4295
$ws.Range("M62").Value = -3671

$ws.Range("H65").Value = 5154.6
$ws.Range("I65").Value = 4295
$ws.Range("K65").Value = 21475
$ws.Range("M65").Value = -18355

$ws.Range("H76").Value = 4998.5
$ws.Range("I76").Value = 4998.5
$ws.Range("K76").Value = 4998.5
$ws.Range("M76").Value = -4683.5

$ws.Range("H79").Value = 4998.5
$ws.Range("I79").Value = 4998.5
$ws.Range("K79").Value = 4998.5
$ws.Range("M79").Value = -3906.5

$ws.Range("H116").Value = 9399.6
$ws.Range("I116").Value = 9249.5
$ws.Range("J116").Value = 9499.666999999999
$ws.Range("K116").Value = 9249.5
$ws.Range("L116").Value = 9499.666999999999
$ws.Range("M116").Value = -5807.5
$ws.Range("N116").Value = -16383.667

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 66664.86
$ws.Range("J24").Value = 66664.86
$ws.Range("L24").Value = 66664.86
$ws.Range("N24").Value = -67412.86

$ws.Range("H45").Value = 1856
$ws.Range("I45").Value = 1748.25
$ws.Range("K45").Value = 1748.25
$ws.Range("M45").Value = -1371.25

$ws.Range("H100").Value = 66664.86
$ws.Range("J100").Value = 66664.86
$ws.Range("L100").Value = 66664.86
$ws.Range("N100").Value = -68828.86

$ws.Range("H110").Value = 1002.75
$ws.Range("I110").Value = 1002.75
$ws.Range("K110").Value = 1002.75
$ws.Range("M110").Value = 1042.25

$ws.Range("H122").Value = 3210.889
$ws.Range("I122").Value = 2974.6
$ws.Range("K122").Value = 8923.799999999999
$ws.Range("M122").Value = -6473.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4672.6665
$ws.Range("I105").Value = 4179
$ws.Range("J105").Value = 5166.3335
$ws.Range("K105").Value = 4179
$ws.Range("L105").Value = 5166.3335
$ws.Range("M105").Value = -2432
$ws.Range("N105").Value = -8660.333500000001

$ws.Range("H134").Value = 988.9167
$ws.Range("I134").Value = 892.0952
$ws.Range("J134").Value = 1666.6666
$ws.Range("K134").Value = 2676.2856
$ws.Range("L134").Value = 4999.9998
$ws.Range("M134").Value = -141.2856000000002
$ws.Range("N134").Value = -10069.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1506
$ws.Range("J16").Value = 2013
$ws.Range("L16").Value = 2013
$ws.Range("N16").Value = -2587

$ws.Range("H31").Value = 2397.6667
$ws.Range("I31").Value = 2243.5
$ws.Range("K31").Value = 2243.5
$ws.Range("M31").Value = -1948.5

$ws.Range("H34").Value = 2397.6667
$ws.Range("I34").Value = 2243.5
$ws.Range("K34").Value = 2243.5
$ws.Range("M34").Value = -2041.5

$ws.Range("H86").Value = 9961487
$ws.Range("I86").Value = 11620735
$ws.Range("K86").Value = 11620735
$ws.Range("M86").Value = -11619612

$ws.Range("H89").Value = 9961487
$ws.Range("I89").Value = 11620735
$ws.Range("K89").Value = 58103675
$ws.Range("M89").Value = -58098059

$ws.Range("H99").Value = 6400.143
$ws.Range("I99").Value = 5800.3335
$ws.Range("K99").Value = 5800.3335
$ws.Range("M99").Value = -4302.3335

$ws.Range("H113").Value = 1506
$ws.Range("J113").Value = 2013
$ws.Range("L113").Value = 2013
$ws.Range("N113").Value = -6353

$ws.Range("H126").Value = 6400.143
$ws.Range("I126").Value = 5800.3335
$ws.Range("K126").Value = 17401.0005
$ws.Range("M126").Value = -14931.0005

$ws.Range("H132").Value = 1975.2903
$ws.Range("I132").Value = 2164.8518
$ws.Range("J132").Value = 695.75
$ws.Range("K132").Value = 6494.555399999999
$ws.Range("L132").Value = 2087.25
$ws.Range("M132").Value = -3964.555399999999
$ws.Range("N132").Value = -7147.25

$ws.Range("H134").Value = 2377.75
$ws.Range("I134").Value = 1737
$ws.Range("K134").Value = 5211
$ws.Range("M134").Value = -2676

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 144.85715
$ws.Range("I6").Value = 152
$ws.Range("K6").Value = 456
$ws.Range("M6").Value = -343

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = ""

$ws.Range("H92").Value = 639.875
$ws.Range("I92").Value = 292.33334
$ws.Range("J92").Value = 848.4
$ws.Range("K92").Value = 877.0000200000001
$ws.Range("L92").Value = 2545.2
$ws.Range("M92").Value = 370.9999799999999
$ws.Range("N92").Value = -5041.2

$ws.Range("H113").Value = 199.6
$ws.Range("I113").Value = 165.33333
$ws.Range("K113").Value = 495.99999
$ws.Range("M113").Value = 1674.00001

$ws.Range("H128").Value = 278225.2
$ws.Range("I128").Value = 278225.2
$ws.Range("K128").Value = 834675.6000000001
$ws.Range("M128").Value = -829695.6000000001

$ws.Range("H131").Value = 1031.2188
$ws.Range("J131").Value = 1029
$ws.Range("L131").Value = 3087
$ws.Range("N131").Value = -13167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2670.8
$ws.Range("J80").Value = 3377.5
$ws.Range("L80").Value = 3377.5
$ws.Range("N80").Value = -5373.5

$ws.Range("H83").Value = 2670.8
$ws.Range("J83").Value = 3377.5
$ws.Range("L83").Value = 16887.5
$ws.Range("N83").Value = -26871.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1999
$ws.Range("I7").Value = 1999
$ws.Range("K7").Value = 1999
$ws.Range("M7").Value = -1887

$ws.Range("H122").Value = 5050
$ws.Range("I122").Value = 4832.636
$ws.Range("K122").Value = 14497.908
$ws.Range("M122").Value = -12047.908

$ws.Range("H126").Value = 1999
$ws.Range("I126").Value = 1999
$ws.Range("K126").Value = 5997
$ws.Range("M126").Value = -3527

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = ""

$ws.Range("H122").Value = 2685
$ws.Range("I122").Value = 2504.111
$ws.Range("J122").Value = 3499
$ws.Range("K122").Value = 7512.333
$ws.Range("L122").Value = 10497
$ws.Range("M122").Value = -5062.333
$ws.Range("N122").Value = -15397

$ws.Range("H132").Value = 2208.238
$ws.Range("I132").Value = 2151.2632
$ws.Range("J132").Value = 2749.5
$ws.Range("K132").Value = 6453.7896
$ws.Range("L132").Value = 8248.5
$ws.Range("M132").Value = -3923.7896
$ws.Range("N132").Value = -13308.5
